# Rename the inline-picture shapes (Pearson logo in both footers, BTec logo
# in the secondary header) back to their "other" generated name, as in the
# commit: the two Pearson-logo copies go image2.png -> image1.png, and the
# BTec logo goes image1.jpg -> image2.jpg.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footers: Pearson logo, "image2.png" -> "image1.png" -------------------
for ($f = 1; $f -le 3; $f++) {
    $ftr = $sec.Footers.Item($f)
    if ($ftr.Exists -and $ftr.Range.InlineShapes.Count -gt 0) {
        for ($i = 1; $i -le $ftr.Range.InlineShapes.Count; $i++) {
            $shp = $ftr.Range.InlineShapes.Item($i)
            if ($shp.AlternativeText -like "*PearsonLogo*") {
                $shp.Name = "image1.png"
            }
        }
    }
}

# --- Headers: BTec logo, "image1.jpg" -> "image2.jpg" ----------------------
for ($h = 1; $h -le 3; $h++) {
    $hdr = $sec.Headers.Item($h)
    if ($hdr.Exists -and $hdr.Range.InlineShapes.Count -gt 0) {
        for ($i = 1; $i -le $hdr.Range.InlineShapes.Count; $i++) {
            $shp = $hdr.Range.InlineShapes.Item($i)
            if ($shp.AlternativeText -like "*BTec_Logo*") {
                $shp.Name = "image2.jpg"
            }
        }
    }
}
